# Groupes2 sample workbook update for EERV Morges:
# - Replace the placeholder paroisse id (9040000000) with the correct one
#   (2010000000) in column D for all data rows (D2:D171).
# - Move the active selection to E8 (matches the author's last cursor
#   position when the workbook was re-saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groupes2")

# Bulk-update the "id paroisse" column with the corrected value.
$ws.Range("D2:D171").Value = 2010000000

# Restore the selected cell/active cell as it was left in the saved file.
$ws.Range("E8").Select()
